# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4)
# - Status columns re-autofit (narrower) to match the shorter replacement text

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Best ColumnWidth value this host can store so the resulting column is the
# autofit width for the new, shorter status text (engine quantizes stored
# XML width to 1/6-character steps, so 12.5 is the nearest input that lands
# on the narrowed column width produced by autofitting the "In Translation" text).
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
